# Regenerated localization-status report: the row that used to show
# 80aa2139-65ec-4332-891f-182726765ae0 ahead of
# 769031a1-6495-46d2-b4a9-ff0fe70fa663 is now re-sorted so that
# 769031a1-... comes first (row 3) and 80aa2139-... comes second (row 4),
# with 769031a1-...'s status updated to "In Translation" on every sheet.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")
$ov.Cells.Item(3, 1).Value = "769031a1-6495-46d2-b4a9-ff0fe70fa663.md"
$ov.Cells.Item(3, 2).Value = "In Translation"
$ov.Cells.Item(3, 3).Value = "In Translation"
$ov.Cells.Item(4, 1).Value = "80aa2139-65ec-4332-891f-182726765ae0.md"
$ov.Cells.Item(4, 2).Value = "In Translation"
$ov.Cells.Item(4, 3).Value = "In Translation"

$i = 0
foreach ($hl in $ov.Hyperlinks) {
    $i++
    if ($i -eq 2) { $hl.TextToDisplay = "769031a1-6495-46d2-b4a9-ff0fe70fa663.md" }
    if ($i -eq 3) { $hl.TextToDisplay = "80aa2139-65ec-4332-891f-182726765ae0.md" }
}

# --- zh-cn detail sheet ----------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Cells.Item(3, 1).Value = "769031a1-6495-46d2-b4a9-ff0fe70fa663.md"
$zh.Cells.Item(3, 2).Value = "In Translation"
$zh.Cells.Item(3, 3).Value = "769031a1-6495-46d2-b4a9-ff0fe70fa663.33aa9fa8dca48c60685a31e77d0147f4dda20079.zh-cn.xlf"
$zh.Cells.Item(3, 4).Value = "2016-01-25 05:58:17"
$zh.Cells.Item(4, 1).Value = "80aa2139-65ec-4332-891f-182726765ae0.md"
$zh.Cells.Item(4, 2).Value = "In Translation"
$zh.Cells.Item(4, 3).Value = "80aa2139-65ec-4332-891f-182726765ae0.baeaf4a7d2ad01b92e5483bb4f31592143d3cc1f.zh-cn.xlf"
$zh.Cells.Item(4, 4).Value = "2016-01-25 05:56:48"

$i = 0
foreach ($hl in $zh.Hyperlinks) {
    $i++
    if ($i -eq 3) { $hl.TextToDisplay = "769031a1-6495-46d2-b4a9-ff0fe70fa663.md" }
    if ($i -eq 4) { $hl.TextToDisplay = "769031a1-6495-46d2-b4a9-ff0fe70fa663.33aa9fa8dca48c60685a31e77d0147f4dda20079.zh-cn.xlf" }
    if ($i -eq 5) { $hl.TextToDisplay = "80aa2139-65ec-4332-891f-182726765ae0.md" }
    if ($i -eq 6) { $hl.TextToDisplay = "80aa2139-65ec-4332-891f-182726765ae0.baeaf4a7d2ad01b92e5483bb4f31592143d3cc1f.zh-cn.xlf" }
}

# --- de-de detail sheet ------------------------------------------------
$de = $wb.Worksheets.Item("de-de")
$de.Cells.Item(3, 1).Value = "769031a1-6495-46d2-b4a9-ff0fe70fa663.md"
$de.Cells.Item(3, 2).Value = "In Translation"
$de.Cells.Item(3, 3).Value = "769031a1-6495-46d2-b4a9-ff0fe70fa663.33aa9fa8dca48c60685a31e77d0147f4dda20079.de-de.xlf"
$de.Cells.Item(3, 4).Value = "2016-01-25 05:58:35"
$de.Cells.Item(4, 1).Value = "80aa2139-65ec-4332-891f-182726765ae0.md"
$de.Cells.Item(4, 2).Value = "In Translation"
$de.Cells.Item(4, 3).Value = "80aa2139-65ec-4332-891f-182726765ae0.baeaf4a7d2ad01b92e5483bb4f31592143d3cc1f.de-de.xlf"
$de.Cells.Item(4, 4).Value = "2016-01-25 05:57:01"

$i = 0
foreach ($hl in $de.Hyperlinks) {
    $i++
    if ($i -eq 3) { $hl.TextToDisplay = "769031a1-6495-46d2-b4a9-ff0fe70fa663.md" }
    if ($i -eq 4) { $hl.TextToDisplay = "769031a1-6495-46d2-b4a9-ff0fe70fa663.33aa9fa8dca48c60685a31e77d0147f4dda20079.de-de.xlf" }
    if ($i -eq 5) { $hl.TextToDisplay = "80aa2139-65ec-4332-891f-182726765ae0.md" }
    if ($i -eq 6) { $hl.TextToDisplay = "80aa2139-65ec-4332-891f-182726765ae0.baeaf4a7d2ad01b92e5483bb4f31592143d3cc1f.de-de.xlf" }
}
